$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6087903704106736
$ws.Range("D2").Value = 0.02425809787929012
$ws.Range("E2").Value = 0.2712026110627281
$ws.Range("F2").Value = 0.4923401975686161
$ws.Range("G2").Value = 0.3349290900893251
$ws.Range("H2").Value = 0.5030912612737737
$ws.Range("I2").Value = 0.8089377786349843
$ws.Range("K2").Value = 0.4053179972980843
$ws.Range("L2").Value = 0.1484486294874188
$ws.Range("N2").Value = 1.44642649302849
$ws.Range("O2").Value = 1.609763374942588
$ws.Range("B3").Value = 0.5823244448273499
$ws.Range("D3").Value = 0.02162209466069243
$ws.Range("E3").Value = 0.2734940017564877
$ws.Range("F3").Value = 0.4897859352744902
$ws.Range("G3").Value = 0.3337919161885381
$ws.Range("H3").Value = 0.5055945222935563
$ws.Range("I3").Value = 0.8173451792023272
$ws.Range("K3").Value = 0.3587545892451374
$ws.Range("L3").Value = 0.1376601091031091
$ws.Range("N3").Value = 1.459389425326044
$ws.Range("O3").Value = 1.612371466496313
$ws.Range("B4").Value = 0.5662903058936877
$ws.Range("D4").Value = 0.01999325439678046
$ws.Range("E4").Value = 0.2750101254944184
$ws.Range("F4").Value = 0.4885247925719582
$ws.Range("G4").Value = 0.3333323892596951
$ws.Range("H4").Value = 0.5073460923495503
$ws.Range("I4").Value = 0.8228703013028849
$ws.Range("K4").Value = 0.3300014714490374
$ws.Range("L4").Value = 0.1310909253713675
$ws.Range("N4").Value = 1.467817406551063
$ws.Range("O4").Value = 1.61492014730257
$ws.Range("B5").Value = 0.5598112063353256
$ws.Range("D5").Value = 0.01932692705285888
$ws.Range("E5").Value = 0.2756554542710319
$ws.Range("F5").Value = 0.4880881794065743
$ws.Range("G5").Value = 0.3332051409330106
$ws.Range("H5").Value = 0.5081138863676244
$ws.Range("I5").Value = 0.8252131388640116
$ws.Range("K5").Value = 0.3182441201999211
$ws.Range("L5").Value = 0.1284279061255376
$ws.Range("N5").Value = 1.471369811275203
$ws.Range("O5").Value = 1.61619706465865
$ws.Range("B6").Value = 0.5587386932683671
$ws.Range("D6").Value = 0.01921613018944157
$ws.Range("E6").Value = 0.2757642723931113
$ws.Range("F6").Value = 0.4880203516329402
$ws.Range("G6").Value = 0.3331876357327488
$ws.Range("H6").Value = 0.5082446421981288
$ws.Range("I6").Value = 0.8256076810616353
$ws.Range("K6").Value = 0.3162894141505603
$ws.Range("L6").Value = 0.1279865621650771
$ws.Range("N6").Value = 1.471966810635138
$ws.Range("O6").Value = 1.616423492111579
$ws.Range("B7").Value = 0.5662027032159642
$ws.Range("D7").Value = 0.01998427840067052
$ws.Range("E7").Value = 0.2750187172443361
$ws.Range("F7").Value = 0.4885185911270469
$ws.Range("G7").Value = 0.3333304301703919
$ws.Range("H7").Value = 0.5073562283004094
$ws.Range("I7").Value = 0.8229015278832801
$ws.Range("K7").Value = 0.3298430695442391
$ws.Range("L7").Value = 0.1310549541822326
$ws.Range("N7").Value = 1.467864837907761
$ws.Range("O7").Value = 1.614936403216063
$ws.Range("B8").Value = 0.5996204390036439
$ws.Range("D8").Value = 0.02335136648915181
$ws.Range("E8").Value = 0.2719700400918921
$ws.Range("F8").Value = 0.4913957708488681
$ws.Range("G8").Value = 0.3344874418772719
$ws.Range("H8").Value = 0.5039098914660585
$ws.Range("I8").Value = 0.8117613460062891
$ws.Range("K8").Value = 0.3892972095574123
$ws.Range("L8").Value = 0.1447174127102642
$ws.Range("N8").Value = 1.450798883610894
$ws.Range("O8").Value = 1.610466090151291
$ws.Range("B9").Value = 0.6668433314977449
$ws.Range("D9").Value = 0.02987108184174048
$ws.Range("E9").Value = 0.2668564139322971
$ws.Range("F9").Value = 0.499473271704737
$ws.Range("G9").Value = 0.3386515123174121
$ws.Range("H9").Value = 0.4988514987966326
$ws.Range("I9").Value = 0.7927933546581691
$ws.Range("K9").Value = 0.5045648711824811
$ws.Range("L9").Value = 0.1719410443612475
$ws.Range("N9").Value = 1.42104762033053
$ws.Range("O9").Value = 1.609213631192489
$ws.Range("B10").Value = 0.7172366265880044
$ws.Range("D10").Value = 0.03460921359829428
$ws.Range("E10").Value = 0.2636244104807446
$ws.Range("F10").Value = 0.5068912506864294
$ws.Range("G10").Value = 0.3428688113565812
$ws.Range("H10").Value = 0.496168240668112
$ws.Range("I10").Value = 0.7806091686689385
$ws.Range("K10").Value = 0.5884156876797419
$ws.Range("L10").Value = 0.1922008393700736
$ws.Range("N10").Value = 1.401448479316297
$ws.Range("O10").Value = 1.612873372878653
$ws.Range("B11").Value = 0.7403747557505653
$ws.Range("D11").Value = 0.03675320403405635
$ws.Range("E11").Value = 0.2622676081305197
$ws.Range("F11").Value = 0.510587818097612
$ws.Range("G11").Value = 0.3450394185883425
$ws.Range("H11").Value = 0.4951712528838925
$ws.Range("I11").Value = 0.7754460072367451
$ws.Range("K11").Value = 0.6263739450358798
$ws.Range("L11").Value = 0.2014728488566391
$ws.Range("N11").Value = 1.393021685364367
$ws.Range("O11").Value = 1.615532793809109
$ws.Range("B12").Value = 0.7491667714591301
$ws.Range("D12").Value = 0.03756340625528765
$ws.Range("E12").Value = 0.2617701001543029
$ws.Range("F12").Value = 0.5120338803817646
$ws.Range("G12").Value = 0.3458976558689812
$ws.Range("H12").Value = 0.4948258221867974
$ws.Range("I12").Value = 0.7735453840331807
$ws.Range("K12").Value = 0.6407203299727939
$ws.Range("L12").Value = 0.2049918101505739
$ws.Range("N12").Value = 1.389900918570188
$ws.Range("O12").Value = 1.616682813923262
$ws.Range("B13").Value = 0.7472719250000353
$ws.Range("D13").Value = 0.03738899008374119
$ws.Range("E13").Value = 0.2618765236679437
$ws.Range("F13").Value = 0.5117203889390396
$ws.Range("G13").Value = 0.3457112054739611
$ws.Range("H13").Value = 0.4948987897344495
$ws.Range("I13").Value = 0.7739522911934635
$ws.Range("K13").Value = 0.6376318196162742
$ws.Range("L13").Value = 0.2042335922030247
$ws.Range("N13").Value = 1.390569906454115
$ws.Range("O13").Value = 1.616428779612477
$ws.Range("B14").Value = 0.7410974805389117
$ws.Range("D14").Value = 0.03681989376572403
$ws.Range("E14").Value = 0.2622263516771834
$ws.Range("F14").Value = 0.5107058599155678
$ws.Range("G14").Value = 0.3451092991708293
$ws.Range("H14").Value = 0.4951421909208022
$ws.Range("I14").Value = 0.7752885485686072
$ws.Range("K14").Value = 0.6275547884455648
$ws.Range("L14").Value = 0.2017621994688028
$ws.Range("N14").Value = 1.392763529727826
$ws.Range("O14").Value = 1.615624542004355
$ws.Range("B15").Value = 0.737319356429964
$ws.Range("D15").Value = 0.0364710858134174
$ws.Range("E15").Value = 0.2624427511335057
$ws.Range("F15").Value = 0.5100904531055903
$ws.Range("G15").Value = 0.3447453388405535
$ws.Range("H15").Value = 0.495295460897708
$ws.Range("I15").Value = 0.7761141488681638
$ws.Range("K15").Value = 0.6213786984590399
$ws.Range("L15").Value = 0.2002494178976804
$ws.Range("N15").Value = 1.39411633852384
$ws.Range("O15").Value = 1.615150538248116
$ws.Range("B16").Value = 0.7157287445969303
$ws.Range("D16").Value = 0.0344688657996528
$ws.Range("E16").Value = 0.2637153610318954
$ws.Range("F16").Value = 0.5066561471833282
$ws.Range("G16").Value = 0.3427320323108347
$ws.Range("H16").Value = 0.4962378919073132
$ws.Range("I16").Value = 0.780954228301642
$ws.Range("K16").Value = 0.5859312127562362
$ws.Range("L16").Value = 0.191595999377796
$ws.Range("N16").Value = 1.402009026180092
$ws.Range("O16").Value = 1.612719581570502
$ws.Range("B17").Value = 0.7025379362240471
$ws.Range("D17").Value = 0.03323761820869464
$ws.Range("E17").Value = 0.2645251019898502
$ws.Range("F17").Value = 0.5046317620470901
$ws.Range("G17").Value = 0.3415615257496967
$ws.Range("H17").Value = 0.4968732851908726
$ws.Range("I17").Value = 0.7840206445851408
$ws.Range("K17").Value = 0.5641371094612566
$ws.Range("L17").Value = 0.1863015662100622
$ws.Range("N17").Value = 1.406976142438648
$ws.Range("O17").Value = 1.611482960195218
$ws.Range("B18").Value = 0.6949711111381305
$ws.Range("D18").Value = 0.03252836563583372
$ws.Range("E18").Value = 0.2650015240651644
$ws.Range("F18").Value = 0.5034977099507074
$ws.Range("G18").Value = 0.3409120126074185
$ws.Range("H18").Value = 0.4972597997578561
$ws.Range("I18").Value = 0.7858200860724018
$ws.Range("K18").Value = 0.5515842708073819
$ws.Range("L18").Value = 0.1832616032644552
$ws.Range("N18").Value = 1.409879120835612
$ws.Range("O18").Value = 1.610865307488751
$ws.Range("B19").Value = 0.6924126000393187
$ws.Range("D19").Value = 0.03228804207407165
$ws.Range("E19").Value = 0.265164667713492
$ws.Range("F19").Value = 0.5031189489308971
$ws.Range("G19").Value = 0.3406961739048882
$ws.Range("H19").Value = 0.4973942844391246
$ws.Range("I19").Value = 0.7864354814523615
$ws.Range("K19").Value = 0.5473311247453978
$ws.Range("L19").Value = 0.1822332314822575
$ws.Range("N19").Value = 1.410869925259021
$ws.Range("O19").Value = 1.61067226202826
$ws.Range("B20").Value = 0.7039400366203097
$ws.Range("D20").Value = 0.03336879780225388
$ws.Range("E20").Value = 0.264437798514205
$ws.Range("F20").Value = 0.504844123682112
$ws.Range("G20").Value = 0.3416836719362522
$ws.Range("H20").Value = 0.4968034679462363
$ws.Range("I20").Value = 0.7836905225702466
$ws.Range("K20").Value = 0.5664589408128791
$ws.Range("L20").Value = 0.1868646249515109
$ws.Range("N20").Value = 1.406442621321649
$ws.Range("O20").Value = 1.611604911879454
$ws.Range("B21").Value = 0.7429102526760119
$ws.Range("D21").Value = 0.03698709707360592
$ws.Range("E21").Value = 0.2621231570222289
$ws.Range("F21").Value = 0.5110025969210312
$ws.Range("G21").Value = 0.3452851089963787
$ws.Range("H21").Value = 0.4950698271354099
$ws.Range("I21").Value = 0.7748945769317857
$ws.Range("K21").Value = 0.6305154121585019
$ws.Range("L21").Value = 0.2024878955942597
$ws.Range("N21").Value = 1.392117302353032
$ws.Range("O21").Value = 1.615856886855823
$ws.Range("B22").Value = 0.7685547296796926
$ws.Range("D22").Value = 0.03934205141376879
$ws.Range("E22").Value = 0.2607053042413536
$ws.Range("F22").Value = 0.5152970791138145
$ws.Range("G22").Value = 0.3478503050110504
$ws.Range("H22").Value = 0.4941239122320411
$ws.Range("I22").Value = 0.769463893024632
$ws.Range("K22").Value = 0.6722188473808899
$ws.Range("L22").Value = 0.2127443013282999
$ws.Range("N22").Value = 1.383164493448014
$ws.Range("O22").Value = 1.619469018122544
$ws.Range("B23").Value = 0.7548519780801826
$ws.Range("D23").Value = 0.03808607907699013
$ws.Range("E23").Value = 0.2614533657886522
$ws.Range("F23").Value = 0.5129803895100977
$ws.Range("G23").Value = 0.3464618582517716
$ws.Range("H23").Value = 0.4946116606109285
$ws.Range("I23").Value = 0.7723332616648264
$ws.Range("K23").Value = 0.6499759676378289
$ws.Range("L23").Value = 0.2072661368644901
$ws.Range("N23").Value = 1.38790530681058
$ws.Range("O23").Value = 1.617464938202858
$ws.Range("B24").Value = 0.7033060948181742
$ws.Range("D24").Value = 0.03330949584064058
$ws.Range("E24").Value = 0.2644772344889201
$ws.Range("F24").Value = 0.50474802217731
$ws.Range("G24").Value = 0.3416283766749757
$ws.Range("H24").Value = 0.4968349662247391
$ws.Range("I24").Value = 0.7838396569824937
$ws.Range("K24").Value = 0.5654093128054853
$ws.Range("L24").Value = 0.1866100540063087
$ws.Range("N24").Value = 1.406683678827175
$ws.Range("O24").Value = 1.611549486905176
$ws.Range("B25").Value = 0.6484794445627529
$ws.Range("D25").Value = 0.02811635181723204
$ws.Range("E25").Value = 0.2681474321122757
$ws.Range("F25").Value = 0.4970274935614114
$ws.Range("G25").Value = 0.3373218562404929
$ws.Range("H25").Value = 0.5000382707451934
$ws.Range("I25").Value = 0.797616895916633
$ws.Range("K25").Value = 0.4735265750697977
$ws.Range("L25").Value = 0.1645305896938538
$ws.Range("N25").Value = 1.428699001520567
$ws.Range("O25").Value = 1.609213631192489
